# Updates cryptocurrency price/volume snapshot values on the "cryptos" sheet.
# Mirrors the GitHub Actions scheduled refresh that re-fetches live market
# data and writes updated Price (column D) and Volume(1h) (column E) cells.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "20.553.40"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +1.57%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.472.46"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +2.23%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.19%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.9618"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +5.37%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "277.23"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +0.66%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3588"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -1.05%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3078"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +0.02%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "1.087"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +6.13%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "39.40"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +1.33%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.06626"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +1.90%  "
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +0.30%  "
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +3.78%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.461"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +2.25%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.167"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +1.98%  "
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +3.29%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001022"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +1.30%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "1.472.34"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +2.36%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.05976"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +6.44%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "68.94"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +1.89%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.483"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +1.44%  "
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +2.34%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "11.24"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +3.72%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.262"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +1.01%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "20.550.57"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +1.57%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "144.30"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +4.02%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.103"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -1.14%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "17.11"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +1.37%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.633.69"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +2.69%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "113.68"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +3.46%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.880"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -0.25%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.07993"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +4.46%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.931"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +2.13%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.8031"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -0.35%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.247"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +10.52%  "
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -0.37%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.05776"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -1.06%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.707"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +0.94%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.02046"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +2.93%  "
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +4.01%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "10.36"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +1.86%  "
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +1.41%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "7.370"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +2.58%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.5257"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +1.05%  "
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +0.84%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "12.11"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +2.68%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "119.01"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +2.16%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.5193"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +2.27%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.806"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +4.23%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.06442"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +1.03%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.9918"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +0.38%  "
